$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 357, pushing the existing rows 357-362 down to 359-364.
$ws.Rows.Item(357).Insert()
$ws.Rows.Item(357).Insert()

# New row 357: Cebolla, Sin especificar, 1a (cosecha)
$ws.Cells.Item(357, 1).Value = 11
$ws.Cells.Item(357, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(357, 3).Value = "Bíobío"
$ws.Cells.Item(357, 4).Value = 44595
$ws.Cells.Item(357, 5).Value = 8
$ws.Cells.Item(357, 6).Value = 100112004
$ws.Cells.Item(357, 7).Value = "Cebolla"
$ws.Cells.Item(357, 8).Value = "Sin especificar"
$ws.Cells.Item(357, 9).Value = "1a (cosecha)"
$ws.Cells.Item(357, 10).Value = 5000
$ws.Cells.Item(357, 11).Value = 1800
$ws.Cells.Item(357, 12).Value = 1800
$ws.Cells.Item(357, 13).Value = 1800
$ws.Cells.Item(357, 14).Value = "`$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item(357, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(357, 16).Value = 90
$ws.Cells.Item(357, 17).Value = 20
$ws.Cells.Item(357, 18).Value = "Hortaliza"

# New row 358: Cebolla, Sin especificar, 2a (cosecha)
$ws.Cells.Item(358, 1).Value = 11
$ws.Cells.Item(358, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(358, 3).Value = "Bíobío"
$ws.Cells.Item(358, 4).Value = 44595
$ws.Cells.Item(358, 5).Value = 8
$ws.Cells.Item(358, 6).Value = 100112004
$ws.Cells.Item(358, 7).Value = "Cebolla"
$ws.Cells.Item(358, 8).Value = "Sin especificar"
$ws.Cells.Item(358, 9).Value = "2a (cosecha)"
$ws.Cells.Item(358, 10).Value = 5000
$ws.Cells.Item(358, 11).Value = 1400
$ws.Cells.Item(358, 12).Value = 1400
$ws.Cells.Item(358, 13).Value = 1400
$ws.Cells.Item(358, 14).Value = "`$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item(358, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(358, 16).Value = 70
$ws.Cells.Item(358, 17).Value = 20
$ws.Cells.Item(358, 18).Value = "Hortaliza"
